$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 1128.2727
$ws.Range("I41").Value = 1783.1111
$ws.Range("J41").Value = 674.9231
$ws.Range("K41").Value = 1783.1111
$ws.Range("L41").Value = 674.9231
$ws.Range("M41").Value = -1343.1111
$ws.Range("N41").Value = -1554.9231

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 44454.793
$ws.Range("J64").Value = 3054.182
$ws.Range("L64").Value = 3054.182
$ws.Range("N64").Value = -3550.182

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 44454.793
$ws.Range("J67").Value = 3054.182
$ws.Range("L67").Value = 3054.182
$ws.Range("N67").Value = -4770.182

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 1937.6666
$ws.Range("I86").Value = 1931.579
$ws.Range("J86").Value = 1952.125
$ws.Range("K86").Value = 1931.579
$ws.Range("L86").Value = 1952.125
$ws.Range("M86").Value = -808.579
$ws.Range("N86").Value = -4198.125

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 1937.6666
$ws.Range("I89").Value = 1931.579
$ws.Range("J89").Value = 1952.125
$ws.Range("K89").Value = 9657.895
$ws.Range("L89").Value = 9760.625
$ws.Range("M89").Value = -4041.895
$ws.Range("N89").Value = -20992.625

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H93").Value = 32533.666
$ws.Range("J93").Value = 32533.666
$ws.Range("L93").Value = 32533.666
$ws.Range("N93").Value = -37525.666

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 168667.5
$ws.Range("I113").Value = 502002.5
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 502002.5
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = -498748.5
$ws.Range("N113").Value = -8508

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H123").Value = 36799.668
$ws.Range("J123").Value = 36799.668
$ws.Range("L123").Value = 36799.668
$ws.Range("N123").Value = -46599.668

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 2612.2698
$ws.Range("J129").Value = 1134.5294
$ws.Range("L129").Value = 3403.5882
$ws.Range("N129").Value = -13403.5882

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 1397.6
$ws.Range("I135").Value = 578.78125
$ws.Range("J135").Value = 3413.1538
$ws.Range("K135").Value = 5209.03125
$ws.Range("L135").Value = 30718.3842
$ws.Range("M135").Value = -2674.03125
$ws.Range("N135").Value = -35788.3842

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 29812.059
$ws.Range("I32").Value = 8156.3125
$ws.Range("K32").Value = 8156.3125
$ws.Range("M32").Value = -7869.3125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2099.1667
$ws.Range("I61").Value = 1673.5555
$ws.Range("K61").Value = 1673.5555
$ws.Range("M61").Value = -1461.5555

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 786.1429000000001
$ws.Range("I74").Value = 691.15
$ws.Range("J74").Value = 1023.625
$ws.Range("K74").Value = 691.15
$ws.Range("L74").Value = 1023.625
$ws.Range("M74").Value = 182.85
$ws.Range("N74").Value = -2771.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 786.1429000000001
$ws.Range("I77").Value = 691.15
$ws.Range("J77").Value = 1023.625
$ws.Range("K77").Value = 3455.75
$ws.Range("L77").Value = 5118.125
$ws.Range("M77").Value = 912.25
$ws.Range("N77").Value = -13854.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1954.129
$ws.Range("J122").Value = 3279.7144
$ws.Range("L122").Value = 9839.143199999999
$ws.Range("N122").Value = -14739.1432

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H123").Value = 48000
$ws.Range("J123").Value = 48000
$ws.Range("L123").Value = 48000
$ws.Range("N123").Value = -57800

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 16329.439
$ws.Range("I132").Value = 19631.219
$ws.Range("J132").Value = 4589.778
$ws.Range("K132").Value = 58893.65700000001
$ws.Range("L132").Value = 13769.334
$ws.Range("M132").Value = -56363.65700000001
$ws.Range("N132").Value = -18829.334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2099.1667
$ws.Range("I136").Value = 1673.5555
$ws.Range("K136").Value = 5020.666499999999
$ws.Range("M136").Value = -2470.666499999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 254277.88
$ws.Range("I105").Value = 201998
$ws.Range("J105").Value = 341411
$ws.Range("K105").Value = 201998
$ws.Range("L105").Value = 341411
$ws.Range("M105").Value = -200251
$ws.Range("N105").Value = -344905

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 14355.581
$ws.Range("I134").Value = 15212.59
$ws.Range("J134").Value = 5999.75
$ws.Range("K134").Value = 45637.77
$ws.Range("L134").Value = 17999.25
$ws.Range("M134").Value = -43102.77
$ws.Range("N134").Value = -23069.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 32314.217
$ws.Range("I31").Value = 560.9697
$ws.Range("J31").Value = 112918.62
$ws.Range("K31").Value = 560.9697
$ws.Range("L31").Value = 112918.62
$ws.Range("M31").Value = -265.9697
$ws.Range("N31").Value = -113508.62

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 32314.217
$ws.Range("I34").Value = 560.9697
$ws.Range("J34").Value = 112918.62
$ws.Range("K34").Value = 560.9697
$ws.Range("L34").Value = 112918.62
$ws.Range("M34").Value = -358.9697
$ws.Range("N34").Value = -113322.62

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2600
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 2600
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 807.125
$ws.Range("I107").Value = 949.2105
$ws.Range("J107").Value = 599.46155
$ws.Range("K107").Value = 949.2105
$ws.Range("L107").Value = 599.46155
$ws.Range("M107").Value = 970.7895
$ws.Range("N107").Value = -4439.46155

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2812.3928
$ws.Range("I132").Value = 3079.8823
$ws.Range("J132").Value = 2399
$ws.Range("K132").Value = 9239.6469
$ws.Range("L132").Value = 7197
$ws.Range("M132").Value = -6709.6469
$ws.Range("N132").Value = -12257

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1218.8235
$ws.Range("J34").Value = 1293.75
$ws.Range("L34").Value = 3881.25
$ws.Range("N34").Value = -4049.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 1650
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 1650
$ws.Range("K58").Value = 0
$ws.Range("L58").ClearContents()
$ws.Range("M58").Value = 4950
$ws.Range("N58").Value = -5206

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 484.6154
$ws.Range("I114").Value = 458.77777
$ws.Range("J114").Value = 542.75
$ws.Range("K114").Value = 1376.33331
$ws.Range("L114").Value = 1628.25
$ws.Range("M114").Value = 1877.66669
$ws.Range("N114").Value = -8136.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 741417.5
$ws.Range("I131").Value = 569
$ws.Range("J131").Value = 881923.25
$ws.Range("K131").Value = 1707
$ws.Range("L131").Value = 2645769.75
$ws.Range("M131").Value = 3333
$ws.Range("N131").Value = -2655849.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2083.8572
$ws.Range("I122").Value = 1709.1177
$ws.Range("J122").Value = 3676.5
$ws.Range("K122").Value = 5127.3531
$ws.Range("L122").Value = 11029.5
$ws.Range("M122").Value = -2677.3531
$ws.Range("N122").Value = -15929.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 7699.8335
$ws.Range("I17").Value = 400
$ws.Range("J17").Value = 14999.667
$ws.Range("K17").Value = 400
$ws.Range("L17").Value = 14999.667
$ws.Range("M17").Value = -230
$ws.Range("N17").Value = -15339.667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H128").Value = 52990
$ws.Range("J128").Value = 52990
$ws.Range("L128").Value = 52990
$ws.Range("N128").Value = -62950

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3651.087
$ws.Range("I132").Value = 4038.75
$ws.Range("K132").Value = 12116.25
$ws.Range("M132").Value = -9586.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H56").Value = 23207.5
$ws.Range("I56").Value = 3950
$ws.Range("J56").Value = 42465
$ws.Range("K56").Value = 3950
$ws.Range("L56").Value = 42465
$ws.Range("M56").Value = -3236
$ws.Range("N56").Value = -43893
